$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.240.52"
$ws.Range("E2").Value = "  +2.27%  "

$ws.Range("D3").Value = "3.387.83"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "585.49"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").Value = "180.14"
$ws.Range("E6").Value = "  +2.61%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("D9").Value = "0.197"
$ws.Range("E9").Value = "  +7.93%  "

$ws.Range("D10").Value = "0.592"
$ws.Range("E10").Value = "  +2.10%  "

$ws.Range("D11").Value = "48.54"
$ws.Range("E11").Value = "  +3.71%  "

$ws.Range("E12").Value = "  +3.97%  "

$ws.Range("D13").Value = "680.15"
$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("D14").Value = "8.66"
$ws.Range("E14").Value = "  +3.07%  "

$ws.Range("D15").Value = "3.926.25"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "69.320.79"
$ws.Range("E16").Value = "  +2.35%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.399.14"
$ws.Range("E17").Value = "  +2.19%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.120"
$ws.Range("E18").Value = "  +1.37%  "

$ws.Range("D19").Value = "17.72"
$ws.Range("E19").Value = "  +1.00%  "

$ws.Range("D20").Value = "11.28"
$ws.Range("E20").Value = "  +2.15%  "

$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("E22").Value = "  -2.46%  "

$ws.Range("D23").Value = "17.19"
$ws.Range("E23").Value = "  +2.08%  "

$ws.Range("D24").Value = "103.11"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  +0.41%  "

$ws.Range("D26").Value = "2.73"
$ws.Range("E26").Value = "  +1.76%  "

$ws.Range("D27").Value = "9.67"
$ws.Range("E27").Value = "  +2.92%  "

$ws.Range("E28").Value = "  +2.65%  "

$ws.Range("D29").Value = "8.78"
$ws.Range("E29").Value = "  +2.91%  "

$ws.Range("D30").Value = "6.94"
$ws.Range("E30").Value = "  -1.22%  "

$ws.Range("D31").Value = "11.14"
$ws.Range("E31").Value = "  +1.19%  "

$ws.Range("D32").Value = "557.30"
$ws.Range("E32").Value = "  -2.29%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "3.58"
$ws.Range("E34").Value = "  +9.94%  "

$ws.Range("E35").Value = "  +1.89%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "3.665.95"
$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("E38").Value = "  +5.36%  "

$ws.Range("D39").Value = "35.60"
$ws.Range("E39").Value = "  +1.51%  "

$ws.Range("D40").Value = "0.0₃0718"
$ws.Range("E40").Value = "  +7.30%  "

$ws.Range("D41").Value = "3.26"
$ws.Range("E41").Value = "  +3.58%  "

$ws.Range("D42").Value = "2.69"
$ws.Range("E42").Value = "  +2.91%  "

$ws.Range("D43").Value = "0.338"
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("D44").Value = "0.0425"
$ws.Range("E44").Value = "  +4.56%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "2.68"
$ws.Range("E46").Value = "  +1.47%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.130"
$ws.Range("E47").Value = "  +1.15%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "1.39"
$ws.Range("E48").Value = "  +4.85%  "

$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "134.21"
$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "2.66"
$ws.Range("E51").Value = "  +2.54%  "
